# Update column G ("K") values on Sheet1 rows 2-32 with the newly
# regenerated K values (computed from the underlying data instead of the
# old Strike# figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 5
    3  = 5
    4  = 10
    5  = 2
    6  = 7
    7  = 5
    8  = 2
    9  = 7
    10 = 2
    11 = 2
    12 = 2
    13 = 7
    14 = 2
    15 = 1
    16 = 2
    17 = 6
    18 = 5
    19 = 3
    20 = 3
    21 = 3
    22 = 3
    23 = 3
    24 = 4
    25 = 4
    26 = 7
    27 = 6
    28 = 10
    29 = 5
    30 = 5
    31 = 2
    32 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
